# [WIP] Add tests for Saliva Check-in validates against all missing fields
# CVDLS-201
#
# Adds new "missing required/optional field" test fixture rows to Sheet1:
#  - row 5  (new): missing Tube ID                         (REQUIRED)  -> red
#  - row 6  (was old row 5): missing Accepted/Rejected      (REQUIRED)  -> red
#  - row 7  (new): duplicate Tube ID TestCheckin0003, missing Well Plate Barcode
#  - row 8  (was old row 6): missing Well Plate Barcode     (OPTIONAL)  -> yellow
#  - row 9  (was old row 7): missing Kit Type                (OPTIONAL)  -> yellow
#  - row 10 (was old row 8): unchanged
#  - row 11 (was old row 9): missing Username                (REQUIRED)  -> red
#
# Column F carries a human-readable note of which field each row tests.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$RED = 255        # RGB(255,0,0)
$YELLOW = 65535   # RGB(255,255,0)

# --- Row 5 (new row): Tests empty Tube ID (REQUIRED) ---
$ws.Range("A5").Value = ""
$ws.Range("A5").Interior.Color = $RED
$ws.Range("B5").Value = "Accepted"
$ws.Range("C5").Value = "TESTBC1"
$ws.Range("D5").Value = "Kit Type 2"
$ws.Range("E5").Value = "jok"
$ws.Range("F5").Value = "Tests empty Tube ID (REQUIRED)"

# --- Row 6 (previously row 5, TestCheckin0003): Tests empty Accepted/Rejected (REQUIRED) ---
$ws.Range("A6").Value = "TestCheckin0003"
$ws.Range("B6").Value = ""
$ws.Range("B6").Interior.Color = $RED
$ws.Range("C6").Value = "TESTBC1"
$ws.Range("D6").Value = "Kit Type 2"
$ws.Range("E6").Value = "jok"
$ws.Range("F6").Value = "Tests empty Accepted/Rejected (REQUIRED)"

# --- Row 7 (new row, duplicate Tube ID TestCheckin0003): empty Well Plate Barcode ---
$ws.Range("A7").Value = "TestCheckin0003"
$ws.Range("B7").Value = "Rejected"
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = "Kit Type 2"
$ws.Range("E7").Value = "jok"

# --- Row 8 (previously row 6, TestCheckin0004): Tests empty Well Plate Barcode (OPTIONAL) ---
$ws.Range("A8").Value = "TestCheckin0004"
$ws.Range("B8").Value = "Accepted"
$ws.Range("C8").Value = ""
$ws.Range("C8").Interior.Color = $YELLOW
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = "jok"
$ws.Range("F8").Value = "Tests empty Well Plate Barcode (OPTIONAL)"

# --- Row 9 (previously row 7, TestCheckin0005): Tests empty Kit Type (OPTIONAL) ---
$ws.Range("A9").Value = "TestCheckin0005"
$ws.Range("B9").Value = "Accepted"
$ws.Range("C9").Value = "TESTBC1"
$ws.Range("D9").Value = ""
$ws.Range("D9").Interior.Color = $YELLOW
$ws.Range("E9").Value = "jok"
$ws.Range("F9").Value = "Tests empty Kit Type (OPTIONAL)"

# --- Row 10 (previously row 8, TestCheckin0006): unchanged content, new position ---
$ws.Range("A10").Value = "TestCheckin0006"
$ws.Range("B10").Value = "Rejected"
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = "Kit Type 2"
$ws.Range("E10").Value = "jok"

# --- Row 11 (previously row 9, TestCheckin0007): Tests empty Username (REQUIRED) ---
$ws.Range("A11").Value = "TestCheckin0007"
$ws.Range("B11").Value = "Accepted"
$ws.Range("C11").Value = "TESTBC1"
$ws.Range("D11").Value = "Kit Type 1"
$ws.Range("E11").Value = ""
$ws.Range("E11").Interior.Color = $RED
$ws.Range("F11").Value = "Tests empty Username (REQUIRED)"

# Reflect the author's final cursor position/selection before save.
$ws.Rows("10:10").Select()

Write-Host "Saliva check-in missing-field test rows added."
